# SSU_ChangingPassword.docx wording fixes
# Uses Find.Execute to locate text, then sets the found Range's .Text
# directly (instead of passing a Replacement string into Execute) so
# that straight double-quotes are not auto-converted into curly quotes.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $found = $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false)
    if ($found) {
        $find.Parent.Text = $new
        return $true
    } else {
        Write-Output "NOT FOUND: $old"
        return $false
    }
}

# 1. Fix double space "Terminal  team" -> "Terminal team"
Replace-Text "This document is intended for members of the Terminal  team during design and testing, and it also can be used for writing the user manual applications." "This document is intended for members of the Terminal team during design and testing, and it also can be used for writing the user manual applications."

# 2. Rewrite the "By pressing the link..." sentence
Replace-Text 'By pressing the link "Forgot password?" on the page "Log in", the user selects the option for changing the password.' 'By pressing the "Settings" button of the top right Profile menu the user is redirected to the password reset page of the website.'

# 3. Rewrite the "A 4-digit code..." sentence
Replace-Text "A 4-digit code is sent to the user’s e-mail address saved in the application database, and it has to be entered in the field marked for it. " "The page has two inputs, one is for the old (current) password and the other input is for the new password that the user wants to change to. "

# 4. Rewrite "Basic user enters correct code. After that..." paragraph (2.2.1 body)
Replace-Text "Basic user enters correct code. After that, he is required to enter a new password that meets the criteria described in the document “5.3.1 Registration”." "Basic user enters everything correctly (password meets the criteria described in the document “5.3.1 Registration”). He is after logged out and redirected to the login page."

# 5. Heading 2.2.2: "...enters incorrect code" -> "...enters incorrect old password"
Replace-Text "2.2.2  Basic user enters incorrect code" "2.2.2  Basic user enters incorrect old password"

# 6. Body under 2.2.2: describe wrong old password instead of wrong code
Replace-Text "Basic user enters a different code from the one sent to his mail. " "Basic user enters a wrong old password, he is notified by an error message about it. "

# 7. Heading 2.2.3: "...enters correct code, but..." -> "...enters correct old password, but..."
Replace-Text "2.2.3  Basic user enters correct code, but the new password does not meet the criteria  " "2.2.3  Basic user enters correct old password, but the new password does not meet the criteria  "

# 8. Body under 2.2.3: correct old password wording + added notification sentence
Replace-Text "Basic user enters correct code, but password doesn’t meet previously described criteria. " "Basic user enters the correct old password, but password doesn’t meet previously described criteria. He is notified by an error message about it. "

# 9. History-of-changes table: fill in the blank "2.1" row
$tbl = $d.Tables(1)
$tbl.Cell(4, 1).Range.Text = "2.1"
$tbl.Cell(4, 2).Range.Text = "19.6.2023."
$tbl.Cell(4, 3).Range.Text = "Fixes in wording"
$tbl.Cell(4, 4).Range.Text = "Jovana Bjelica," + [char]11 + "Andrej Dujovi" + [char]0x107
